$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto table
# with the latest scraped figures. A few Price cells (D12, D23, D47) carry a
# trailing zero that would otherwise be silently dropped if Excel parsed the
# text as a plain number (e.g. "5.20" -> 5.2), so those are entered with a
# leading apostrophe to force them to stay literal text, same as typing them
# into Excel by hand.
$ws.Range('D2').Value = '62.057.80'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.423.13'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '563.33'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '143.79'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('D9').Value = '2.422.26'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '''5.20'
$ws.Range('E12').Value = '  -3.66%  '
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').Value = '2.859.12'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '61.959.45'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '2.410.69'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').Value = '11.28'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = '323.74'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').Value = '6.84'
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '67.19'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('D26').Value = '8.83'
$ws.Range('E26').Value = '  -2.10%  '
$ws.Range('D27').Value = '557.28'
$ws.Range('E27').Value = '  -5.61%  '
$ws.Range('D28').Value = '2.543.36'
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').Value = '8.21'
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('D32').Value = '1.39'
$ws.Range('E32').Value = '  -5.20%  '
$ws.Range('E33').Value = '  -1.88%  '
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D37').Value = '4.76'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D40').Value = '152.33'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  -3.84%  '
$ws.Range('D45').Value = '147.58'
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('D47').Value = '''0.0530'
$ws.Range('E47').Value = '  -1.79%  '
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').Value = '19.91'
$ws.Range('E49').Value = '  -2.50%  '
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('E51').Value = '  -0.62%  '
